$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '88.176.37'
$ws.Range('E2').Value = '  -1.55%  '
$ws.Range('D3').Value = '3.068.94'
$ws.Range('E3').Value = '  -3.78%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '210.04'
$ws.Range('E5').Value = '  -3.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '618.37'
$ws.Range('E6').Value = '  -1.13%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.370'
$ws.Range('E7').Value = '  -4.98%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.802'
$ws.Range('E8').Value = '  +16.24%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('D10').Value = '3.065.46'
$ws.Range('E10').Value = '  -3.74%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.593'
$ws.Range('E11').Value = '  +3.66%  '
$ws.Range('E12').Value = '  -0.41%  '
$ws.Range('E13').Value = '  -8.02%  '
$ws.Range('E14').Value = '  -2.00%  '
$ws.Range('D15').Value = '87.862.49'
$ws.Range('E15').Value = '  -1.57%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '31.93'
$ws.Range('E16').Value = '  -4.63%  '
$ws.Range('D17').Value = '3.630.00'
$ws.Range('E17').Value = '  -3.88%  '
$ws.Range('D18').Value = '3.073.33'
$ws.Range('E18').Value = '  -3.77%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.26'
$ws.Range('E19').Value = '  -5.17%  '
$ws.Range('E20').Value = '  -11.43%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.25'
$ws.Range('E21').Value = '  -1.92%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '420.39'
$ws.Range('E22').Value = '  -3.47%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.12'
$ws.Range('E23').Value = '  -6.00%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.87'
$ws.Range('E24').Value = '  -4.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.46'
$ws.Range('E25').Value = '  +2.36%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.70'
$ws.Range('E26').Value = '  -1.82%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '81.68'
$ws.Range('E27').Value = '  -0.43%  '
$ws.Range('E28').Value = '  -4.79%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('E30').Value = '  +8.64%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.169'
$ws.Range('E31').Value = '  +6.14%  '
$ws.Range('E32').Value = '  -5.31%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '507.35'
$ws.Range('E33').Value = '  -6.95%  '
$ws.Range('E34').Value = '  -12.29%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.75'
$ws.Range('E35').Value = '  -3.48%  '
$ws.Range('E36').Value = '  -6.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.24'
$ws.Range('E37').Value = '  -6.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '22.22'
$ws.Range('E38').Value = '  -0.65%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.131'
$ws.Range('E39').Value = '  +2.09%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '22.23'
$ws.Range('E40').Value = '  -0.68%  '
$ws.Range('E41').Value = '  +0.33%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.359'
$ws.Range('E43').Value = '  -3.96%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '147.44'
$ws.Range('E44').Value = '  +0.72%  '
$ws.Range('E45').Value = '  -6.46%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.133'
$ws.Range('E46').Value = '  +6.63%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '43.33'
$ws.Range('E47').Value = '  -0.95%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0688'
$ws.Range('E48').Value = '  +14.57%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.701'
$ws.Range('E49').Value = '  -7.68%  '
$ws.Range('E50').Value = '  -5.50%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '155.54'
$ws.Range('E51').Value = '  -10.38%  '
